# Auto-generated edit script: update Leve profit-calculation sheets
# per scheduled runner refresh of current market prices.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 17149.75
$ws.Range("I86").Value = 2799.5
$ws.Range("J86").Value = 31500
$ws.Range("K86").Value = 2799.5
$ws.Range("L86").Value = 31500
$ws.Range("M86").Value = -1676.5
$ws.Range("N86").Value = -33746
$ws.Range("H89").Value = 17149.75
$ws.Range("I89").Value = 2799.5
$ws.Range("J89").Value = 31500
$ws.Range("K89").Value = 13997.5
$ws.Range("L89").Value = 157500
$ws.Range("M89").Value = -8381.5
$ws.Range("N89").Value = -168732
$ws.Range("H116").Value = 114077.555
$ws.Range("I116").Value = 202619.6
$ws.Range("J116").Value = 3400
$ws.Range("K116").Value = 202619.6
$ws.Range("L116").Value = 3400
$ws.Range("M116").Value = -199177.6
$ws.Range("N116").Value = -10284
$ws.Range("H137").Value = 1452.6111
$ws.Range("I137").Value = 1370.5834
$ws.Range("J137").Value = 1616.6666
$ws.Range("K137").Value = 4111.7502
$ws.Range("L137").Value = 4849.9998
$ws.Range("M137").Value = -1561.7502
$ws.Range("N137").Value = -9949.9998
$ws.Range("H138").Value = 2161.2307
$ws.Range("I138").Value = 1092.1428
$ws.Range("J138").Value = 2555.1052
$ws.Range("K138").Value = 3276.4284
$ws.Range("L138").Value = 7665.3156
$ws.Range("M138").Value = 1863.5716
$ws.Range("N138").Value = -17945.3156

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3119.125
$ws.Range("I32").Value = 3119.125
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 3119.125
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -2832.125
$ws.Range("H45").Value = 1715.7693
$ws.Range("I45").Value = 1500.5
$ws.Range("J45").Value = 2433.3333
$ws.Range("K45").Value = 1500.5
$ws.Range("L45").Value = 2433.3333
$ws.Range("M45").Value = -1123.5
$ws.Range("N45").Value = -3187.3333
$ws.Range("H55").Value = 22251.7
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 22251.7
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 22251.7
$ws.Range("N55").Value = -22881.7
$ws.Range("M55").ClearContents()
$ws.Range("H61").Value = 3039.1936
$ws.Range("I61").Value = 2101.0557
$ws.Range("J61").Value = 4338.154
$ws.Range("K61").Value = 2101.0557
$ws.Range("L61").Value = 4338.154
$ws.Range("M61").Value = -1889.0557
$ws.Range("N61").Value = -4762.154
$ws.Range("H110").Value = 2248.9
$ws.Range("I110").Value = 1329.8334
$ws.Range("J110").Value = 3627.5
$ws.Range("K110").Value = 1329.8334
$ws.Range("L110").Value = 3627.5
$ws.Range("M110").Value = 715.1666
$ws.Range("N110").Value = -7717.5
$ws.Range("H136").Value = 3039.1936
$ws.Range("I136").Value = 2101.0557
$ws.Range("J136").Value = 4338.154
$ws.Range("K136").Value = 6303.1671
$ws.Range("L136").Value = 13014.462
$ws.Range("M136").Value = -3753.1671
$ws.Range("N136").Value = -18114.462

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H40").Value = 44989
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 44989
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 44989
$ws.Range("N40").Value = -45519
$ws.Range("H99").Value = 658
$ws.Range("I99").Value = 658
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 658
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 840
$ws.Range("H134").Value = 20836834
$ws.Range("I134").Value = 3638
$ws.Range("J134").Value = 41670030
$ws.Range("K134").Value = 10914
$ws.Range("L134").Value = 125010090
$ws.Range("M134").Value = -8379
$ws.Range("N134").Value = -125015160

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H80").Value = 23593
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 23593
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 23593
$ws.Range("N80").Value = -25839
$ws.Range("H83").Value = 23593
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 23593
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 70779
$ws.Range("N83").Value = -82011
$ws.Range("H111").Value = 54999
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 54999
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 54999
$ws.Range("N111").Value = -63179
$ws.Range("H122").Value = 1586.4706
$ws.Range("I122").Value = 1640.6154
$ws.Range("J122").Value = 1410.5
$ws.Range("K122").Value = 4921.8462
$ws.Range("L122").Value = 4231.5
$ws.Range("M122").Value = -2471.8462
$ws.Range("N122").Value = -9131.5
$ws.Range("H132").Value = 5035.5303
$ws.Range("I132").Value = 5035.915
$ws.Range("J132").Value = 5032.2856
$ws.Range("K132").Value = 15107.745
$ws.Range("L132").Value = 15096.8568
$ws.Range("M132").Value = -12577.745
$ws.Range("N132").Value = -20156.8568
$ws.Range("H134").Value = 3450175.5
$ws.Range("I134").Value = 1723.6
$ws.Range("J134").Value = 25003000
$ws.Range("K134").Value = 5170.799999999999
$ws.Range("L134").Value = 75009000
$ws.Range("M134").Value = -2635.799999999999
$ws.Range("N134").Value = -75014070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 768.6
$ws.Range("I5").Value = 646.9
$ws.Range("J5").Value = 1012
$ws.Range("K5").Value = 1940.7
$ws.Range("L5").Value = 3036
$ws.Range("M5").Value = -1828.7
$ws.Range("N5").Value = -3260
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H113").Value = 554.8889
$ws.Range("I113").Value = 312.5
$ws.Range("J113").Value = 748.8
$ws.Range("K113").Value = 937.5
$ws.Range("L113").Value = 2246.4
$ws.Range("M113").Value = 1232.5
$ws.Range("N113").Value = -6586.4
$ws.Range("H135").Value = 768.6
$ws.Range("I135").Value = 646.9
$ws.Range("J135").Value = 1012
$ws.Range("K135").Value = 5822.099999999999
$ws.Range("L135").Value = 9108
$ws.Range("M135").Value = -3287.099999999999
$ws.Range("N135").Value = -14178

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2957.6667
$ws.Range("I122").Value = 2988.05
$ws.Range("J122").Value = 2805.75
$ws.Range("K122").Value = 8964.150000000001
$ws.Range("L122").Value = 8417.25
$ws.Range("M122").Value = -6514.150000000001
$ws.Range("N122").Value = -13317.25
$ws.Range("H126").Value = 6364.273
$ws.Range("I126").Value = 16499
$ws.Range("J126").Value = 4112.1113
$ws.Range("K126").Value = 49497
$ws.Range("L126").Value = 12336.3339
$ws.Range("M126").Value = -47027
$ws.Range("N126").Value = -17276.3339

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4636.875
$ws.Range("I40").Value = 4300.857
$ws.Range("J40").Value = 6989
$ws.Range("K40").Value = 4300.857
$ws.Range("L40").Value = 6989
$ws.Range("M40").Value = -4164.857
$ws.Range("N40").Value = -7261
$ws.Range("H122").Value = 3381.875
$ws.Range("I122").Value = 3254
$ws.Range("J122").Value = 3595
$ws.Range("K122").Value = 9762
$ws.Range("L122").Value = 10785
$ws.Range("M122").Value = -7312
$ws.Range("N122").Value = -15685
$ws.Range("H136").Value = 64586530
$ws.Range("I136").Value = 3653.4285
$ws.Range("J136").Value = 516666660
$ws.Range("K136").Value = 10960.2855
$ws.Range("L136").Value = 1549999980
$ws.Range("M136").Value = -8410.2855
$ws.Range("N136").Value = -1550005080

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 1875.3077
$ws.Range("I14").Value = 1870.8182
$ws.Range("J14").Value = 1900
$ws.Range("K14").Value = 1870.8182
$ws.Range("L14").Value = 1900
$ws.Range("M14").Value = -1702.8182
$ws.Range("N14").Value = -2236
$ws.Range("H33").Value = 11500
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 11500
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 11500
$ws.Range("N33").Value = -12000
$ws.Range("H36").Value = 11500
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 11500
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 11500
$ws.Range("N36").Value = -12000
$ws.Range("H41").Value = 100397.664
$ws.Range("I41").Value = 32374
$ws.Range("J41").Value = 108900.625
$ws.Range("K41").Value = 32374
$ws.Range("L41").Value = 108900.625
$ws.Range("M41").Value = -31984
$ws.Range("N41").Value = -109680.625
$ws.Range("H113").Value = 840.16
$ws.Range("I113").Value = 898
$ws.Range("J113").Value = 657
$ws.Range("K113").Value = 2694
$ws.Range("L113").Value = 1971
$ws.Range("M113").Value = -524
$ws.Range("N113").Value = -6311
$ws.Range("H136").Value = 1590.4117
$ws.Range("I136").Value = 1448.8148
$ws.Range("J136").Value = 2136.5715
$ws.Range("K136").Value = 4346.4444
$ws.Range("L136").Value = 6409.7145
$ws.Range("M136").Value = -1796.4444
$ws.Range("N136").Value = -11509.7145
